$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163, shifting existing rows 163-215 down to 164-216.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with a new weekly observation,
# matching the template of the surrounding rows for this market/category.
$ws.Range("A163").Value = 3
$ws.Range("B163").Value = "Femacal de La Calera"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 44524
$ws.Range("E163").Value = 5
$ws.Range("F163").Value = 100112039
$ws.Range("G163").Value = "Ciboulette"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 160
$ws.Range("K163").Value = 1500
$ws.Range("L163").Value = 1500
$ws.Range("M163").Value = 1500
$ws.Range("N163").Value = "$/docena de atados"
$ws.Range("O163").Value = "Provincia de Quillota"
$ws.Range("P163").Value = 500
$ws.Range("Q163").Value = 3
$ws.Range("R163").Value = "Hortaliza"

# Match the date number format used by the rest of the Fecha column.
$ws.Range("D163").NumberFormat = "YYYY-MM-DD HH:MM:SS"
